$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{A="ECs"; D="Inflammatory-Mac"; E=2; F=1; G=19.4786585; H=38.957317; I=0.01644248566400343; J=0.01108359890151296; K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846; O=0.1293260700537641; P=0.1293260700537641; Q=3.167528544863667; R=19.005171269182; S=0.002126442052840921; T=0.001433398287984888},
    @{A="ECs"; D="Neutrophils"; E=2; F=1; G=19.4786585; H=38.957317; I=0.01644248566400343; J=0.01108359890151296; K=3; L=1; M=0.8767803333333334; N=2.630341; O=0.6972931302732585; P=0.6972931302732585; Q=17.07850469251617; R=102.471028155097; S=0.01146523229812613; T=0.007728517372729221},
    @{A="ECs"; D="Resolving-Mac"; E=2; F=1; G=19.4786585; H=38.957317; I=0.01644248566400343; J=0.01108359890151296; K=3; L=1; M=0.21801; N=0.65403; O=0.1733807996729775; P=0.1733807996729775; Q=4.246542339585001; R=25.47925403751; S=0.002850811313036383; T=0.001921683240798852},
    @{A="FAPs"; D="Inflammatory-Mac"; E=3; F=1; G=392.0055033333333; H=1176.01651; I=0.3309029145291901; J=0.3345840089140918; K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846; O=0.1293260700537641; P=0.1293260700537641; Q=63.74610559305111; R=573.71495033746; S=0.04279437350539676; T=0.04327043497569309},
    @{A="FAPs"; D="Neutrophils"; E=3; F=1; G=392.0055033333333; H=1176.01651; I=0.3309029145291901; J=0.3345840089140918; K=3; L=1; M=0.8767803333333334; N=2.630341; O=0.6972931302732585; P=0.6972931302732585; Q=343.7027158811011; R=3093.32444292991; S=0.2307363290886034; T=0.2333031309150829},
    @{A="FAPs"; D="Resolving-Mac"; E=3; F=1; G=392.0055033333333; H=1176.01651; I=0.3309029145291901; J=0.3345840089140918; K=3; L=1; M=0.21801; N=0.65403; O=0.1733807996729775; P=0.1733807996729775; Q=85.46111978169999; R=769.1500780353; S=0.05737221193518988; T=0.05801044302331586},
    @{A="Inflammatory-Mac"; D="Inflammatory-Mac"; E=3; F=1; G=375.1018676666667; H=1125.305603; I=0.3166340783504202; J=0.3201564405802684; K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846; O=0.1293260700537641; P=0.1293260700537641; Q=60.99731524457089; R=548.975837201138; S=0.04094904099815549; T=0.04140457426264757},
    @{A="Inflammatory-Mac"; D="Neutrophils"; E=3; F=1; G=375.1018676666667; H=1125.305603; I=0.3166340783504202; J=0.3201564405802684; K=3; L=1; M=0.8767803333333334; N=2.630341; O=0.6972931302732585; P=0.6972931302732585; Q=328.8819405667359; R=2959.937465100623; S=0.2207867676441527; T=0.2232428866293598},
    @{A="Inflammatory-Mac"; D="Resolving-Mac"; E=3; F=1; G=375.1018676666667; H=1125.305603; I=0.3166340783504202; J=0.3201564405802684; K=3; L=1; M=0.21801; N=0.65403; O=0.1733807996729775; P=0.1733807996729775; Q=81.77595817001001; R=735.98362353009; S=0.05489826970811205; T=0.05550897968826102},
    @{A="MuSCs"; D="Inflammatory-Mac"; E=2; F=1; G=19.6220475; H=39.244095; I=0.01656352436781744; J=0.01116518902553968; K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846; O=0.1293260700537641; P=0.1293260700537641; Q=3.190845794895; R=19.14507476937; S=0.002142095512729588; T=0.001443950018080463},
    @{A="MuSCs"; D="Neutrophils"; E=2; F=1; G=19.6220475; H=39.244095; I=0.01656352436781744; J=0.01116518902553968; K=3; L=1; M=0.8767803333333334; N=2.630341; O=0.6972931302732585; P=0.6972931302732585; Q=17.2042253477325; R=103.225352086395; S=0.01154963175479282; T=0.007785409605711192},
    @{A="MuSCs"; D="Resolving-Mac"; E=2; F=1; G=19.6220475; H=39.244095; I=0.01656352436781744; J=0.01116518902553968; K=3; L=1; M=0.21801; N=0.65403; O=0.1733807996729775; P=0.1733807996729775; Q=4.277802575475; R=25.66681545285; S=0.002871797100295037; T=0.001935829401748021},
    @{A="Neutrophils"; D="Inflammatory-Mac"; E=3; F=1; G=133.4172743333333; H=400.2518229999999; I=0.1126212886044614; J=0.1138741321875775; K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846; O=0.1293260700537641; P=0.1293260700537641; Q=21.69569453813978; R=195.261250843258; S=0.01456486865960577; T=0.01472689399660224},
    @{A="Neutrophils"; D="Neutrophils"; E=3; F=1; G=133.4172743333333; H=400.2518229999999; I=0.1126212886044614; J=0.1138741321875775; K=3; L=1; M=0.8767803333333334; N=2.630341; O=0.6972931302732585; P=0.6972931302732585; Q=116.9776422624048; R=1052.798780361643; S=0.07853005086641297; T=0.07940365009022671},
    @{A="Neutrophils"; D="Resolving-Mac"; E=3; F=1; G=133.4172743333333; H=400.2518229999999; I=0.1126212886044614; J=0.1138741321875775; K=3; L=1; M=0.21801; N=0.65403; O=0.1733807996729775; P=0.1733807996729775; Q=29.08629997741; R=261.77669979669; S=0.01952636907844271; T=0.01974358810074852},
    @{A="Resolving-Mac"; D="Inflammatory-Mac"; E=3; F=1; G=245.0287756666667; H=735.086327; I=0.2068357084841073; J=0.2091366303910096; K=2; L=0.6666666666666666; M=0.1626153333333333; N=0.487846; O=0.1293260700537641; P=0.1293260700537641; Q=39.84543603129355; R=358.608924281642; S=0.0267492493250356; T=0.02704681851275589},
    @{A="Resolving-Mac"; D="Neutrophils"; E=3; F=1; G=245.0287756666667; H=735.086327; I=0.2068357084841073; J=0.2091366303910096; K=3; L=1; M=0.8767803333333334; N=2.630341; O=0.6972931302732585; P=0.6972931302732585; Q=214.8364116052786; R=1933.527704447507; S=0.1442251186211704; T=0.1458295356601486},
    @{A="Resolving-Mac"; D="Resolving-Mac"; E=3; F=1; G=245.0287756666667; H=735.086327; I=0.2068357084841073; J=0.2091366303910096; K=3; L=1; M=0.21801; N=0.65403; O=0.1733807996729775; P=0.1733807996729775; Q=53.41872338309; R=480.76851044781; S=0.03586134053790137; T=0.03626027621810517},
)

$rowNum = 2
foreach ($r in $rows) {
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = "Ccl2"
    $ws.Cells.Item($rowNum, 3).Value = "Ccr3"
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $r.R
    $ws.Cells.Item($rowNum, 19).Value = $r.S
    $ws.Cells.Item($rowNum, 20).Value = $r.T
    $rowNum++
}
